# Correct Random field length
# - "Random (16 byte)" fields are actually 4 bytes -> rename to "Random (4 byte)"
# - Row 12 (Sensor Data packet) is missing a "Length (2 byte)" field, which needs
#   to be inserted before "NodeID (2 byte)", shifting the remaining fields right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 12: insert "Length (2 byte)" field before "NodeID (2 byte)" ---
# Before: C12=NodeID (2 byte)  D12=Sensor Data (n)  E12=Random (16 byte)  F12=CRC (4 byte)
# After:  C12=Length (2 byte)  D12=NodeID (2 byte)  E12=Random (4 byte)  F12=Sensor Data (n)  G12=CRC (4 byte)
$ws.Range("G12").Value = "CRC (4 byte)"
$ws.Range("F12").Value = "Sensor Data (n)"
$ws.Range("C12").Value = "Length (2 byte)"
$ws.Range("E12").Value = "Random (4 byte)"
$ws.Range("D12").Value = "NodeID (2 byte)"

# --- Simple rename of "Random (16 byte)" -> "Random (4 byte)" ---
$ws.Range("C8").Value = "Random (4 byte)"
$ws.Range("I9").Value = "Random (4 byte)"
$ws.Range("J15").Value = "Random (4 byte)"

# --- Sheet view adjustments to match the saved state ---
$ws.Range("I10").Select()
